# Added unity commands to main script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B (Question) and column C (ConditionType) values per row
$values = @{
    2  = @{ B = 6;  C = 2 }
    3  = @{ B = 9;  C = 1 }
    4  = @{ B = 5;  C = 1 }
    5  = @{ B = 7;  C = 1 }
    6  = @{ B = 3;  C = 2 }
    7  = @{ B = 8;  C = 2 }
    8  = @{ B = 1;  C = 1 }
    9  = @{ B = 10; C = 2 }
    10 = @{ B = 2;  C = 1 }
    11 = @{ B = 4;  C = 2 }
    12 = @{ B = 8;  C = 2 }
    13 = @{ B = 10; C = 1 }
    14 = @{ B = 2;  C = 1 }
    15 = @{ B = 1;  C = 1 }
    16 = @{ B = 6;  C = 1 }
    17 = @{ B = 3;  C = 1 }
    18 = @{ B = 4;  C = 1 }
    19 = @{ B = 5;  C = 1 }
    20 = @{ B = 9;  C = 1 }
    21 = @{ B = 7;  C = 1 }
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 2).Value = $values[$row].B
    $ws.Cells.Item($row, 3).Value = $values[$row].C
}

# Update the selected range to reflect the new selection in the saved view
$ws.Range("A1:C11").Select()
